# Update LDLC prices history: add a new snapshot column P, mirroring the
# latest prior snapshot column O (header timestamp + values), matching the
# "Update LDLC prices history" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 204
$srcCol = 15   # O
$dstCol = 16   # P

# Header cell (P1): same style as the other header cells (copy format from
# O1, which already carries the bold/border/centered header style), then set
# the new snapshot timestamp.
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Cells.Item(1, $dstCol).Value = "2026-01-28 07:18:05"

# Data rows: mirror column O's value into column P for every product row.
for ($r = 2; $r -le $lastRow; $r++) {
    $srcCell = $ws.Cells.Item($r, $srcCol)
    $dstCell = $ws.Cells.Item($r, $dstCol)
    $v = $srcCell.Value()
    if ($v -eq $null -or $v -eq "") {
        # Rows with no price yet (O is blank) stay blank in P too, but the
        # cell must still exist so the sheet's used range grows to column P.
        $dstCell.Style = "Normal"
    } else {
        $dstCell.Value = $v
    }
}
